# Assignments.xlsx - "Added link to spreadsheet"
# Adds a third assignment row (Assignment 3 / MovieAPI) and turns the
# GitHub-link cells (B2:B4) into real hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3's submission date moved one day later --------------------------
$ws.Range("C3").Value = 45181

# --- New row 4: copy row 3's formatting down, then fill in the content ----
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Rows(4).RowHeight = 35.25

$ws.Range("A4").Value = "Assignment 3"
$ws.Range("B4").Value = "https://github.com/Vasanth30e/Assignment_Phase_4/tree/master/Assignment_3/MovieAPI"
$ws.Range("C4").Value = 45182

# --- Turn the GitHub-link cells into real hyperlinks -----------------------
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/Vasanth30e/Assignment_Phase_4/tree/master/Assignment_3/MovieAPI")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/Vasanth30e/Assignment_Phase_4/tree/master/Assingnment_2/UserAuthenticationTesting")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/Vasanth30e/Assignment_Phase_4/tree/master/Assignment_1/LoginLib")

# --- Match the saved selection state ---------------------------------------
$ws.Range("B9").Select()
